# atualizei dados da bibi e add
#
# Refresh the "vendas_atipicas" anomaly figures (estoque_atualizado /
# media_vendas / desvio_padrao on several already-present rows — their
# id_venda stays the same, only the computed stats drifted) and append the
# newly detected atypical sales for 2025-08-11 (two more BEMOL S/A rows for
# product 10130, sale ids 396518 and 396572).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be stored as text even when the string looks like a
    # number/date (e.g. "2025-08-11", "396518") — Excel would otherwise
    # auto-coerce it to a date serial / numeric literal on plain assignment.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 3 (2025-07-29 / BEMOL S/A / 389699) — estoque_atualizado drifted
$ws.Range("G3").Value = -71

# Row 4 (2025-07-29 / BEMOL S/A / 390273)
$ws.Range("G4").Value = -77

# Row 5 (2025-07-30 / BEMOL S/A / 390286)
$ws.Range("G5").Value = -27
$ws.Range("I5").Value = 0.26

# Row 6 (2025-07-30 / BEMOL S/A / 390878)
$ws.Range("G6").Value = -77
$ws.Range("I6").Value = 0.15

# Row 7 (2025-07-31 / BEMOL S/A / 391921)
$ws.Range("G7").Value = -134
$ws.Range("H7").Value = 1.26
$ws.Range("I7").Value = 1.43

# Row 8 (2025-08-04 / BEMOL S/A / 393760)
$ws.Range("G8").Value = -1282
$ws.Range("H8").Value = 1.06
$ws.Range("I8").Value = 0.27

# Row 10 (2025-08-07 / BEMOL S/A / 394429)
$ws.Range("G10").Value = -23
$ws.Range("H10").Value = 1.23
$ws.Range("I10").Value = 0.53

# Row 11 (2025-08-08 / BEMOL S/A / 394429)
$ws.Range("G11").Value = -1282
$ws.Range("H11").Value = 1.06
$ws.Range("I11").Value = 0.27

# New row 12: 2025-08-11 anomaly for sale id 396518
Set-TextValue $ws.Range("A12") "2025-08-11"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "BEMOL S/A"
Set-TextValue $ws.Range("D12") "396518"
$ws.Range("E12").Value = 10130
$ws.Range("F12").Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Range("G12").Value = -1282
$ws.Range("H12").Value = 1.06
$ws.Range("I12").Value = 0.27

# New row 13: 2025-08-11 anomaly for sale id 396572
Set-TextValue $ws.Range("A13") "2025-08-11"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "BEMOL S/A"
Set-TextValue $ws.Range("D13") "396572"
$ws.Range("E13").Value = 10130
$ws.Range("F13").Value = "FONE DE OUVIDO SEM FIO A GOLD V5.3"
$ws.Range("G13").Value = -1282
$ws.Range("H13").Value = 1.06
$ws.Range("I13").Value = 0.27
